$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF") with the same style as the other header cells (copy format from H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 and J1 data columns for rows 2 through 69
$data = @(
    @(3, 5),
    @(4, 4),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(7, 7),
    @(10, 10),
    @(10, 10),
    @(10, 10),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(10, 10),
    @(5, 6),
    @(6, 7),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(5, 6),
    @(6, 6),
    @(5, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(7, 8),
    @(6, 7),
    @(6, 7),
    @(4, 5),
    @(8, 8),
    @(4, 7),
    @(6, 6),
    @(5, 6),
    @(8, 8),
    @(6, 6),
    @(7, 9),
    @(5, 7),
    @(8, 9),
    @(7, 7),
    @(6, 7),
    @(5, 6),
    @(8, 8),
    @(6, 7),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(7, 7)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
